$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Patient name fields (row 6)
$ws.Range("A6").Value = "ASICONA"
$ws.Range("C6").Value = "DEL BARRIO"
$ws.Range("E6").Value = "NAZARIO"
$ws.Range("G6").Value = "BALTAZAR"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "/201761083"

# Birth date, age, place of birth (row 12)
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "1982-07-30"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "35"
$ws.Range("H12").Value = "CHAJUIL QUICHE "

# Civil status, occupation, cedula (row 14)
$ws.Range("A14").Value = "null"
$ws.Range("D14").Value = "DISEÑADOR GRAFICO"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "2645288161405"

# Emergency contact (row 20)
$ws.Range("A20").Value = "MARIA RAIMUNDO "
$ws.Range("F20").Value = "ESPOSA"
$ws.Range("H20").Value = "7AV 7-10 RES LA EURECA"
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "49308448"
